# Update LR-pairs sheet with new TPM-derived values.
# The "ECs" sending-cluster rows are inserted ahead of the existing "FAPs"
# sending-cluster rows, and the FAPs rows themselves get refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Wnt2", "Fzd8", "ECs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 3.083576666666666, 9.250729999999999, 0.2272509363535097, 0.2272509363535097, 0.03300146534555555, 0.2970131881099999, 0.0009607368530725774, 0.0009607368530725774),
    @("ECs", "Wnt2", "Fzd8", "FAPs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 6.453984666666667, 19.361954, 0.4756405360586227, 0.4756405360586227, 0.06907269523088889, 0.6216542570779999, 0.00201084052342853, 0.00201084052342853),
    @("ECs", "Wnt2", "Fzd8", "MuSCs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 4.031477000000001, 12.094431, 0.2971085275878677, 0.2971085275878677, 0.04314621067966667, 0.388315896117, 0.001256070124048959, 0.001256070124048959),
    @("FAPs", "Wnt2", "Fzd8", "ECs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 3.083576666666666, 9.250729999999999, 0.2272509363535097, 0.2272509363535097, 7.773104729946666, 69.95794256951999, 0.2262901995004371, 0.2262901995004371),
    @("FAPs", "Wnt2", "Fzd8", "FAPs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 6.453984666666667, 19.361954, 0.4756405360586227, 0.4756405360586227, 16.26925617961067, 146.423305616496, 0.4736296955351942, 0.4736296955351941),
    @("FAPs", "Wnt2", "Fzd8", "MuSCs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 4.031477000000001, 12.094431, 0.2971085275878677, 0.2971085275878677, 10.162579473416, 91.46321526074402, 0.2958524574638187, 0.2958524574638187)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
